# Apply the "minor changes for the revision" edit to the verification matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row before row 26: "upSampleAndGetMeanExprPreSeqPreseq"
#    (a sibling entry to the existing upSampleAndGetMeanExprPreSeqZTNB row).
#    The blank row naturally inherits the plain s="3"/s="5" formatting from
#    the row above it, which is exactly what this row needs.
# ---------------------------------------------------------------------------
$ws.Rows("26:26").Insert()
$ws.Range("A26").Value = "upSampleAndGetMeanExprPreSeqPreseq"
$ws.Range("B26").Value = "Used in figure 5 - the results look reasonable, no further tests have been done. The code is close to identical to that of upSampleAndGetMeanExprPreSeqZTNB"
$ws.Rows("26:26").RowHeight = 30

# ---------------------------------------------------------------------------
# 2) Insert a new row before (the now-shifted) row 34: "BinomialDownsampling.R"
#    This lands right before the "GenBugSummary.R" row. The inserted blank
#    row inherits formatting from the row above (predPreSeq, s="3") instead
#    of from GenBugSummary.R (s="7"), so copy the exact formatting down from
#    the still-unmoved GenBugSummary.R row instead of toggling individual
#    font/alignment properties (which would otherwise mint extra unused
#    cell styles).
# ---------------------------------------------------------------------------
$ws.Rows("34:34").Insert()
$ws.Range("A35").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A34").Value = "BinomialDownsampling.R"
$ws.Range("B34").Value = "TCR0011"

# ---------------------------------------------------------------------------
# 3) Highlight the figure-generating-code rows in yellow (new fill), to flag
#    them for follow-up. Touch the "horizontal=left" rows (GenFig1.R /
#    GenFig1Data.R) before the plain bold rows (GenFig2_S4_S5.R onward) so
#    that the two distinct new cell styles get allocated in that order.
# ---------------------------------------------------------------------------
$ws.Range("A37:A38").Interior.Color = 65535
$ws.Range("A39:A47").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 4) Append two follow-up note rows (only column A is populated), leaving a
#    gap of two blank rows after the table.
# ---------------------------------------------------------------------------
$ws.Range("A51").Value = "FIX file names and add the new figures here!"
$ws.Range("A52").Value = "Also test the simulated data!"

# Final selection, matching the saved view state.
$ws.Range("A52").Select() | Out-Null
